$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Target table (pairing simulation re-run with one more entry, so the
# whole results table is regenerated, not just a single appended row).
# Columns: A = countdown number, B = From, C = To, D = Drone (all text
# except A, which is numeric).
$data = @(
    @(21, "2", "12", "1"),
    @(20, "12", "2", "1"),
    @(19, "3", "14", "2"),
    @(18, "7", "3", "2"),
    @(17, "14", "7", "2"),
    @(16, "2", "15", "5"),
    @(15, "15", "2", "5"),
    @(14, "2", "16", "6"),
    @(13, "16", "17", "6"),
    @(12, "17", "2", "6"),
    @(11, "2", "18", "9"),
    @(10, "13", "2", "9"),
    @(9, "18", "13", "9"),
    @(8, "5", "6", "10"),
    @(7, "6", "5", "10"),
    @(6, "4", "8", "12"),
    @(5, "8", "4", "12"),
    @(4, "1", "11", "13"),
    @(3, "11", "1", "13"),
    @(2, "4", "9", "14"),
    @(1, "9", "10", "14"),
    @(0, "10", "4", "14")
)

# Every text value 1..18 needed below already exists as a text cell
# somewhere in the current table; remember one source cell per value.
$sourceMap = @{
    "1"  = "B17"
    "2"  = "B5"
    "3"  = "B2"
    "4"  = "D2"
    "5"  = "D5"
    "6"  = "C19"
    "7"  = "B3"
    "8"  = "B9"
    "9"  = "C15"
    "10" = "C8"
    "11" = "C17"
    "12" = "D15"
    "13" = "B6"
    "14" = "C2"
    "15" = "B12"
    "16" = "C11"
    "17" = "C13"
    "18" = "C5"
}

# Stage each of those text values (value only, real text type preserved)
# off to one side before anything in A2:D22 gets overwritten.
$stage = @{}
$stageCol = 6
foreach ($key in $sourceMap.Keys) {
    $ws.Range($sourceMap[$key]).Copy()
    $cell = $ws.Cells.Item(26, $stageCol)
    $cell.PasteSpecial(-4163)
    $stage[$key] = $cell
    $stageCol = $stageCol + 1
}

# Grab the bold/centered/bordered format already used by column A so it
# can be stamped onto every row, including the new ones.
$ws.Range("A2").Copy()
$ws.Range("A2:A23").PasteSpecial(-4122)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]

    for ($col = 2; $col -le 4; $col++) {
        $key = $entry[$col - 1]
        $stage[$key].Copy()
        $ws.Cells.Item($row, $col).PasteSpecial(-4163)
    }

    $row = $row + 1
}

# Remove the staging cells.
$ws.Range("F26:W26").Value = ""
